$wb = $excel.ActiveWorkbook

# --- Sheet ALC (index 1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H69").Value = 3720
$ws.Range("J69").Value = 3293.3333
$ws.Range("L69").Value = 9879.999899999999
$ws.Range("N69").Value = -11627.9999
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H72").Value = 3720
$ws.Range("J72").Value = 3293.3333
$ws.Range("L72").Value = 29639.9997
$ws.Range("N72").Value = -38375.9997
$ws.Range("H75").Value = 32650
$ws.Range("J75").Value = 32650
$ws.Range("L75").Value = 32650
$ws.Range("N75").Value = -34522
$ws.Range("H78").Value = 32650
$ws.Range("J78").Value = 32650
$ws.Range("L78").Value = 97950
$ws.Range("N78").Value = -107310
$ws.Range("H137").Value = 1348.44
$ws.Range("I137").Value = 1368.8422
$ws.Range("J137").Value = 1283.8334
$ws.Range("K137").Value = 4106.5266
$ws.Range("L137").Value = 3851.5002
$ws.Range("M137").Value = -1556.5266
$ws.Range("N137").Value = -8951.5002

# --- Sheet ARM (index 2) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 33208.355
$ws.Range("I2").Value = 981.9666999999999
$ws.Range("J2").Value = 1000000
$ws.Range("K2").Value = 981.9666999999999
$ws.Range("L2").Value = 1000000
$ws.Range("M2").Value = -868.9666999999999
$ws.Range("N2").Value = -1000226
$ws.Range("H32").Value = 5780.41
$ws.Range("I32").Value = 3895.9456
$ws.Range("J32").Value = 27451.75
$ws.Range("K32").Value = 3895.9456
$ws.Range("L32").Value = 27451.75
$ws.Range("M32").Value = -3608.9456
$ws.Range("N32").Value = -28025.75
$ws.Range("H45").Value = 149543.14
$ws.Range("I45").Value = 201360.6
$ws.Range("J45").Value = 19999.5
$ws.Range("K45").Value = 201360.6
$ws.Range("L45").Value = 19999.5
$ws.Range("M45").Value = -200983.6
$ws.Range("N45").Value = -20753.5
$ws.Range("H61").Value = 1783.6666
$ws.Range("I61").Value = 1235.8
$ws.Range("J61").Value = 2982.125
$ws.Range("K61").Value = 1235.8
$ws.Range("L61").Value = 2982.125
$ws.Range("M61").Value = -1023.8
$ws.Range("N61").Value = -3406.125
$ws.Range("H74").Value = 1167.3871
$ws.Range("I74").Value = 1155.8889
$ws.Range("J74").Value = 1245
$ws.Range("K74").Value = 1155.8889
$ws.Range("L74").Value = 1245
$ws.Range("M74").Value = -281.8888999999999
$ws.Range("N74").Value = -2993
$ws.Range("H77").Value = 1167.3871
$ws.Range("I77").Value = 1155.8889
$ws.Range("J77").Value = 1245
$ws.Range("K77").Value = 5779.4445
$ws.Range("L77").Value = 6225
$ws.Range("M77").Value = -1411.4445
$ws.Range("N77").Value = -14961
$ws.Range("H116").Value = 33208.355
$ws.Range("I116").Value = 981.9666999999999
$ws.Range("J116").Value = 1000000
$ws.Range("K116").Value = 981.9666999999999
$ws.Range("L116").Value = 1000000
$ws.Range("M116").Value = 1312.0333
$ws.Range("N116").Value = -1004588
$ws.Range("H122").Value = 1621.5454
$ws.Range("I122").Value = 1553.7894
$ws.Range("J122").Value = 1713.5
$ws.Range("K122").Value = 4661.3682
$ws.Range("L122").Value = 5140.5
$ws.Range("M122").Value = -2211.3682
$ws.Range("N122").Value = -10040.5
$ws.Range("H132").Value = 2515.1191
$ws.Range("I132").Value = 2092.8157
$ws.Range("K132").Value = 6278.4471
$ws.Range("M132").Value = -3748.4471
$ws.Range("H136").Value = 1783.6666
$ws.Range("I136").Value = 1235.8
$ws.Range("J136").Value = 2982.125
$ws.Range("K136").Value = 3707.4
$ws.Range("L136").Value = 8946.375
$ws.Range("M136").Value = -1157.4
$ws.Range("N136").Value = -14046.375

# --- Sheet BSM (index 3) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 33208.355
$ws.Range("I3").Value = 981.9666999999999
$ws.Range("J3").Value = 1000000
$ws.Range("K3").Value = 981.9666999999999
$ws.Range("L3").Value = 1000000
$ws.Range("M3").Value = -867.9666999999999
$ws.Range("N3").Value = -1000228
$ws.Range("H86").Value = 101381.38
$ws.Range("I86").Value = 80051.86
$ws.Range("J86").Value = 144040.42
$ws.Range("K86").Value = 80051.86
$ws.Range("L86").Value = 144040.42
$ws.Range("M86").Value = -78928.86
$ws.Range("N86").Value = -146286.42
$ws.Range("H89").Value = 101381.38
$ws.Range("I89").Value = 80051.86
$ws.Range("J89").Value = 144040.42
$ws.Range("K89").Value = 400259.3
$ws.Range("L89").Value = 720202.1000000001
$ws.Range("M89").Value = -394643.3
$ws.Range("N89").Value = -731434.1000000001
$ws.Range("H108").Value = 39684
$ws.Range("J108").Value = 39684
$ws.Range("L108").Value = 39684
$ws.Range("N108").Value = -47364
$ws.Range("H134").Value = 4307.923
$ws.Range("I134").Value = 3834.8
$ws.Range("J134").Value = 5885
$ws.Range("K134").Value = 11504.4
$ws.Range("L134").Value = 17655
$ws.Range("M134").Value = -8969.400000000001
$ws.Range("N134").Value = -22725

# --- Sheet CRP (index 4) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 17371.637
$ws.Range("J68").Value = 17371.637
$ws.Range("L68").Value = 17371.637
$ws.Range("N68").Value = -18869.637
$ws.Range("H71").Value = 17371.637
$ws.Range("J71").Value = 17371.637
$ws.Range("L71").Value = 52114.91099999999
$ws.Range("N71").Value = -59602.91099999999

# --- Sheet CUL (index 5) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1111.5385
$ws.Range("J34").Value = 1187.5
$ws.Range("L34").Value = 3562.5
$ws.Range("N34").Value = -3730.5
$ws.Range("H127").Value = 1095.091
$ws.Range("J127").Value = 1095.091
$ws.Range("L127").Value = 3285.273
$ws.Range("N127").Value = -13205.273
$ws.Range("H131").Value = 7579.6113
$ws.Range("J131").Value = 7655.7866
$ws.Range("L131").Value = 22967.3598
$ws.Range("N131").Value = -33047.35980000001

# --- Sheet GSM (index 6) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1739.5555
$ws.Range("I126").Value = 1718.9524
$ws.Range("J126").Value = 1811.6666
$ws.Range("K126").Value = 5156.857199999999
$ws.Range("L126").Value = 5434.9998
$ws.Range("M126").Value = -2686.857199999999
$ws.Range("N126").Value = -10374.9998
$ws.Range("H136").Value = 15822.182
$ws.Range("J136").Value = 15822.182
$ws.Range("L136").Value = 47466.546
$ws.Range("N136").Value = -52566.546

# --- Sheet LTW (index 7) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1485.081
$ws.Range("I136").Value = 1310.75
$ws.Range("J136").Value = 2600.8
$ws.Range("K136").Value = 3932.25
$ws.Range("L136").Value = 7802.400000000001
$ws.Range("M136").Value = -1382.25
$ws.Range("N136").Value = -12902.4

# --- Sheet WVR (index 8) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 34915
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 34915
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 34915
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -36787
$ws.Range("H78").Value = 34915
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 34915
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 104745
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -114105
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H136").Value = 1642.3959
$ws.Range("I136").Value = 557.5
$ws.Range("J136").Value = 2560.3845
$ws.Range("K136").Value = 1672.5
$ws.Range("L136").Value = 7681.1535
$ws.Range("M136").Value = 877.5
$ws.Range("N136").Value = -12781.1535
